# Append the HOSubcortical atlas rows (rows 29-42) to Sheet 1, mirroring the
# rows already present for the GSS atlas. Columns are:
#   A = ROI (number), B = Atlas, C = Hemisphere, D = Parcel, E = Selectivity

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(29, 4,  "HOSubcortical", "Left",  "Left Thalamus",     "Exclusive Responsive"),
    @(30, 5,  "HOSubcortical", "Left",  "Left Caudate",      ""),
    @(31, 6,  "HOSubcortical", "Left",  "Left Putamen",      ""),
    @(32, 7,  "HOSubcortical", "Left",  "Left Pallidum",     ""),
    @(33, 9,  "HOSubcortical", "Left",  "Left Hippocampus",  "Selective"),
    @(34, 10, "HOSubcortical", "Left",  "Left Amygdala",     "Selective"),
    @(35, 11, "HOSubcortical", "Left",  "Left Accumbens",    ""),
    @(36, 15, "HOSubcortical", "Right", "Right Thalamus",    ""),
    @(37, 16, "HOSubcortical", "Right", "Right Caudate",     ""),
    @(38, 17, "HOSubcortical", "Right", "Right Putamen",     ""),
    @(39, 18, "HOSubcortical", "Right", "Right Pallidum",    ""),
    @(40, 19, "HOSubcortical", "Right", "Right Hippocampus", "Selective"),
    @(41, 20, "HOSubcortical", "Right", "Right Amygdala",    "Selective"),
    @(42, 21, "HOSubcortical", "Right", "Right Accumbens",   "")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    if ($row[5] -ne "") {
        $ws.Cells.Item($r, 5).Value = $row[5]
    }
}
